$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32: Automata for the People
$ws.Range("H32").Value = 3836.125
$ws.Range("I32").Value = 6398.3335
$ws.Range("J32").Value = 2298.8
$ws.Range("K32").Value = 6398.3335
$ws.Range("L32").Value = 2298.8
$ws.Range("M32").Value = -6072.3335
$ws.Range("N32").Value = -2950.8

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 1915.5834
$ws.Range("I111").Value = 1612.5714
$ws.Range("J111").Value = 2339.8
$ws.Range("K111").Value = 4837.7142
$ws.Range("L111").Value = 7019.400000000001
$ws.Range("M111").Value = -1770.7142
$ws.Range("N111").Value = -13153.4

# Row 116: Growing Up
$ws.Range("H116").Value = 22225422
$ws.Range("I116").Value = 50002000
$ws.Range("J116").Value = 4159.8
$ws.Range("K116").Value = 50002000
$ws.Range("L116").Value = 4159.8
$ws.Range("M116").Value = -49998558
$ws.Range("N116").Value = -11043.8

# Row 135: For Tired Minds
$ws.Range("H135").Value = 745.0571
$ws.Range("I135").Value = 720.1613
$ws.Range("J135").Value = 938
$ws.Range("K135").Value = 6481.4517
$ws.Range("L135").Value = 8442
$ws.Range("M135").Value = -3946.4517
$ws.Range("N135").Value = -13512

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1119.6207
$ws.Range("I137").Value = 1031.4286
$ws.Range("J137").Value = 1351.125
$ws.Range("K137").Value = 3094.2858
$ws.Range("L137").Value = 4053.375
$ws.Range("M137").Value = -544.2857999999997
$ws.Range("N137").Value = -9153.375

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2616.9207
$ws.Range("J138").Value = 4427.7393
$ws.Range("L138").Value = 13283.2179
$ws.Range("N138").Value = -23563.2179

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 13521.46
$ws.Range("I32").Value = 13384.195
$ws.Range("J32").Value = 15100
$ws.Range("K32").Value = 13384.195
$ws.Range("L32").Value = 15100
$ws.Range("M32").Value = -13097.195
$ws.Range("N32").Value = -15674

# Row 121: Shield to Shield
$ws.Range("H121").Value = 27996.666
$ws.Range("J121").Value = 27996.666
$ws.Range("L121").Value = 27996.666
$ws.Range("N121").Value = -31490.666

$ws = $wb.Worksheets.Item("CRP")
# Row 21: Nightmare on My Street
$ws.Range("H21").Value = 4750
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 4750
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 4750
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -5220

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 4833.136
$ws.Range("I22").Value = 6465.8125
$ws.Range("J22").Value = 479.33334
$ws.Range("K22").Value = 6465.8125
$ws.Range("L22").Value = 479.33334
$ws.Range("M22").Value = -6115.8125
$ws.Range("N22").Value = -1179.33334

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1406.8334
$ws.Range("I58").Value = 1378.9231
$ws.Range("J58").Value = 1479.4
$ws.Range("K58").Value = 1378.9231
$ws.Range("L58").Value = 1479.4
$ws.Range("M58").Value = -1175.9231
$ws.Range("N58").Value = -1885.4

# Row 68: Do You Even String Bow
$ws.Range("H68").Value = 31600
$ws.Range("J68").Value = 31600
$ws.Range("L68").Value = 31600
$ws.Range("N68").Value = -33098

# Row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value = 31600
$ws.Range("J71").Value = 31600
$ws.Range("L71").Value = 94800
$ws.Range("N71").Value = -102288

# Row 86: Birch, Please
$ws.Range("H86").Value = 8127.0454
$ws.Range("I86").Value = 10898.083
$ws.Range("J86").Value = 4801.8
$ws.Range("K86").Value = 10898.083
$ws.Range("L86").Value = 4801.8
$ws.Range("M86").Value = -9775.083000000001
$ws.Range("N86").Value = -7047.8

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 8127.0454
$ws.Range("I89").Value = 10898.083
$ws.Range("J89").Value = 4801.8
$ws.Range("K89").Value = 54490.415
$ws.Range("L89").Value = 24009
$ws.Range("M89").Value = -48874.415
$ws.Range("N89").Value = -35241

# Row 99: O Pine
$ws.Range("H99").Value = 2075.4443
$ws.Range("I99").Value = 2151
$ws.Range("J99").Value = 1697.6666
$ws.Range("K99").Value = 2151
$ws.Range("L99").Value = 1697.6666
$ws.Range("M99").Value = -653
$ws.Range("N99").Value = -4693.6666

# Row 126: A Better Conductor
$ws.Range("H126").Value = 2075.4443
$ws.Range("I126").Value = 2151
$ws.Range("J126").Value = 1697.6666
$ws.Range("K126").Value = 6453
$ws.Range("L126").Value = 5092.9998
$ws.Range("M126").Value = -3983
$ws.Range("N126").Value = -10032.9998

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1685.3784
$ws.Range("I132").Value = 1308.3214
$ws.Range("J132").Value = 2858.4443
$ws.Range("K132").Value = 3924.9642
$ws.Range("L132").Value = 8575.332900000001
$ws.Range("M132").Value = -1394.9642
$ws.Range("N132").Value = -13635.3329

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1689.0938
$ws.Range("I134").Value = 1216.1072
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 3648.3216
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -1113.3216
$ws.Range("N134").Value = -20070

# Row 136: Turali Quality
$ws.Range("H136").Value = 1406.8334
$ws.Range("I136").Value = 1378.9231
$ws.Range("J136").Value = 1479.4
$ws.Range("K136").Value = 4136.7693
$ws.Range("L136").Value = 4438.200000000001
$ws.Range("M136").Value = -1586.7693
$ws.Range("N136").Value = -9538.200000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 92: Oh No Udon
$ws.Range("H92").Value = 612.5
$ws.Range("I92").Value = 833
$ws.Range("J92").Value = 502.25
$ws.Range("K92").Value = 2499
$ws.Range("L92").Value = 1506.75
$ws.Range("M92").Value = -1251
$ws.Range("N92").Value = -4002.75

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 605.7
$ws.Range("I113").Value = 556.2222
$ws.Range("J113").Value = 646.1818
$ws.Range("K113").Value = 1668.6666
$ws.Range("L113").Value = 1938.5454
$ws.Range("M113").Value = 501.3334
$ws.Range("N113").Value = -6278.5454

$ws = $wb.Worksheets.Item("GSM")
# Row 3: Needful Rings
$ws.Range("H3").Value = 3377075
$ws.Range("I3").Value = 3859386.2
$ws.Range("J3").Value = 898
$ws.Range("K3").Value = 3859386.2
$ws.Range("L3").Value = 898
$ws.Range("M3").Value = -3859270.2
$ws.Range("N3").Value = -1130

# Row 51: When We Were Blings
$ws.Range("H51").Value = 15183.167
$ws.Range("J51").Value = 15183.167
$ws.Range("L51").Value = 15183.167
$ws.Range("N51").Value = -16201.167

# Row 119: Bulking Up
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676

# Row 123: Workplace Workout
$ws.Range("H123").Value = 12408.667
$ws.Range("J123").Value = 12408.667
$ws.Range("L123").Value = 12408.667
$ws.Range("N123").Value = -17308.667

# Row 132: On Board for Lar
$ws.Range("H132").Value = 2973.7273
$ws.Range("I132").Value = 2714
$ws.Range("J132").Value = 4999.6
$ws.Range("K132").Value = 8142
$ws.Range("L132").Value = 14998.8
$ws.Range("M132").Value = -5612
$ws.Range("N132").Value = -20058.8

# Row 138: Orders Anonymous
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Row 139: Ringing Gratitude
$ws.Range("H139").Value = 650000
$ws.Range("J139").Value = 650000
$ws.Range("L139").Value = 650000
$ws.Range("N139").Value = -660280

# Row 140: The Right Rod
$ws.Range("H140").Value = 107675
$ws.Range("J140").Value = 107675
$ws.Range("L140").Value = 107675
$ws.Range("N140").Value = -118035

# Row 141: Mask Maker
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 1800.6666
$ws.Range("I16").Value = 1713.125
$ws.Range("J16").Value = 2501
$ws.Range("K16").Value = 1713.125
$ws.Range("L16").Value = 2501
$ws.Range("M16").Value = -1543.125
$ws.Range("N16").Value = -2841

# Row 130: Generous Soles
$ws.Range("H130").Value = 24347
$ws.Range("J130").Value = 24347
$ws.Range("L130").Value = 24347
$ws.Range("N130").Value = -34387

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3319.5227
$ws.Range("I132").Value = 2981.5862
$ws.Range("J132").Value = 3972.8667
$ws.Range("K132").Value = 8944.758600000001
$ws.Range("L132").Value = 11918.6001
$ws.Range("M132").Value = -6414.758600000001
$ws.Range("N132").Value = -16978.6001

# Row 133: The Perfect Accessory
$ws.Range("H133").Value = 70442
$ws.Range("J133").Value = 70442
$ws.Range("L133").Value = 70442
$ws.Range("N133").Value = -75502

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 17721090
$ws.Range("I136").Value = 22223624
$ws.Range("J136").Value = 836584.5
$ws.Range("K136").Value = 66670872
$ws.Range("L136").Value = 2509753.5
$ws.Range("M136").Value = -66668322
$ws.Range("N136").Value = -2514853.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2058.4211
$ws.Range("I132").Value = 1712.4706
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5137.4118
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2607.4118
$ws.Range("N132").Value = -20057

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 1650.8695
$ws.Range("I136").Value = 1650.8695
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4952.6085
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2402.6085
$ws.Range("N136").ClearContents()

Write-Host "All updates applied"